# Renames the embedded logo pictures in the document's first-page header
# and the two footers, matching:
#   footer (first page)  : image1.png -> image2.png
#   footer (default)     : image1.png -> image2.png
#   header (first page)  : image2.jpg -> image1.jpg
#
# The pictures live as inline drawings inside the header/footer ranges, so
# they are reached via Sections(1).Headers / .Footers rather than
# ActiveDocument.InlineShapes (which only covers the main body story).

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footers -----------------------------------------------------------
# wdHeaderFooterPrimary = 1 (default footer), wdHeaderFooterFirstPage = 2
$ftrDefault = $sec.Footers.Item(1)
if ($ftrDefault.Exists) {
    $shp = $ftrDefault.Range.InlineShapes.Item(1)
    $shp.Name = "image2.png"
}

$ftrFirst = $sec.Footers.Item(2)
if ($ftrFirst.Exists) {
    $shp = $ftrFirst.Range.InlineShapes.Item(1)
    $shp.Name = "image2.png"
}

# --- Header (first page) ------------------------------------------------
$hdrFirst = $sec.Headers.Item(2)
if ($hdrFirst.Exists) {
    $shp = $hdrFirst.Range.InlineShapes.Item(1)
    $shp.Name = "image1.jpg"
}
